# costing.xlsx -- "FAQs, Embroidery Base Price Calculation,
# Dynamic Embroidery Base Price Calculations"
#
# The Stitch Count input (Input!B4) was updated from 0 to 50000, which
# ripples through the Output sheet's embroidery base-price formulas
# (Fabric Cost / Patch Attach Cost / Total Cost). The active sheet/
# selection also moved from Material -> Input while working on the
# Input sheet, and finally to Output!C19 to review the recalculated
# total.

$wb = $excel.ActiveWorkbook

$wsMaterial = $wb.Worksheets.Item("Material")
$wsInput    = $wb.Worksheets.Item("Input")
$wsOutput   = $wb.Worksheets.Item("Output")

# Dynamic embroidery base price calculation input: Stitch Count 0 -> 50000.
$wsInput.Range("B4").Value = 50000

# Review the recalculated total on the Output sheet.
$wsOutput.Activate()
$wsOutput.Range("C19").Select()

# Finish on the Input sheet with B4 (the cell that was just edited) selected
# and active -- this becomes the workbook's active tab.
$wsInput.Activate()
$wsInput.Range("B4").Select()
